$wb = $excel.ActiveWorkbook

# --- STAT sheet: add two new error code rows ---
$statWs = $wb.Worksheets.Item("STAT")
$statWs.Cells.Item(2, 1).Value = 721
$statWs.Cells.Item(2, 2).Value = "Invalid User id"
$statWs.Cells.Item(3, 1).Value = 722
$statWs.Cells.Item(3, 2).Value = "Invalid gounp number"

# --- TOURNAMENT sheet: select column B (whole column) ---
$tourWs = $wb.Worksheets.Item("TOURNAMENT")
[void]$tourWs.Activate()
[void]$tourWs.Columns.Item(2).Select()

# --- STAT sheet: make it the active tab, with B4 selected ---
[void]$statWs.Activate()
[void]$statWs.Range("B4").Select()
